$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.089.19"
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = "'3.458.75"
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'579.06"
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').Value = "'149.41"
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('D12').Value = "'4.051.78"
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').Value = "'28.63"
$ws.Range('E14').Value = '  -3.84%  '
$ws.Range('D15').Value = "'3.461.12"
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = "'63.133.29"
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').Value = "'6.46"
$ws.Range('E18').Value = '  +2.83%  '
$ws.Range('D19').Value = "'14.48"
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('E20').Value = '  -2.85%  '
$ws.Range('D21').Value = "'388.84"
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').Value = "'0.563"
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('D23').Value = "'74.74"
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = "'3.593.45"
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('E26').Value = '  -3.19%  '
$ws.Range('D27').Value = "'0.184"
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('D28').Value = "'7.69"
$ws.Range('E28').Value = '  -2.25%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  -2.10%  '
$ws.Range('D31').Value = "'2.11"
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').Value = "'23.40"
$ws.Range('E33').Value = '  -1.80%  '
$ws.Range('D34').Value = "'1.33"
$ws.Range('E34').Value = '  -5.97%  '
$ws.Range('D35').Value = "'1.64"
$ws.Range('E35').Value = '  +3.58%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').Value = "'32.02"
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('D39').Value = "'170.24"
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('D40').Value = "'3.495.31"
$ws.Range('D41').Value = "'0.0776"
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('E43').Value = '  +0.87%  '
$ws.Range('D44').Value = "'1.71"
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('E46').Value = '  -2.92%  '
$ws.Range('D47').Value = "'2.587.75"
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('D48').Value = "'2.31"
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('E49').Value = '  +1.96%  '
$ws.Range('D50').Value = "'22.69"
$ws.Range('E50').Value = '  -4.66%  '
$ws.Range('E51').Value = '  +0.02%  '
